$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 117.666664
$ws.Range("I5").Value = 100.25
$ws.Range("K5").Value = 100.25
$ws.Range("M5").Value = 14.75

$ws.Range("H18").Value = 457.14285
$ws.Range("I18").Value = 400
$ws.Range("K18").Value = 400
$ws.Range("M18").Value = -116

$ws.Range("H51").Value = 15555.444
$ws.Range("I51").Value = 17000
$ws.Range("J51").Value = 13749.75
$ws.Range("K51").Value = 17000
$ws.Range("L51").Value = 13749.75
$ws.Range("M51").Value = -16516
$ws.Range("N51").Value = -14717.75

$ws.Range("H101").Value = 803.6
$ws.Range("I101").Value = 803.6
$ws.Range("K101").Value = 2410.8
$ws.Range("M101").Value = -788.8000000000002

$ws.Range("H106").Value = 144407.72
$ws.Range("I106").Value = 144407.72
$ws.Range("K106").Value = 144407.72
$ws.Range("M106").Value = -143776.72

$ws.Range("H113").Value = 7487.9653
$ws.Range("I113").Value = 7936.9585
$ws.Range("K113").Value = 7936.9585
$ws.Range("M113").Value = -4682.9585

$ws.Range("H132").Value = 1977.9773
$ws.Range("I132").Value = 1792.875
$ws.Range("K132").Value = 5378.625
$ws.Range("M132").Value = -2848.625

$ws.Range("H135").Value = 967.63635
$ws.Range("I135").Value = 1045
$ws.Range("J135").Value = 477.66666
$ws.Range("K135").Value = 9405
$ws.Range("L135").Value = 4298.99994
$ws.Range("M135").Value = -6870
$ws.Range("N135").Value = -9368.99994

$ws.Range("H138").Value = 4245.18
$ws.Range("J138").Value = 4937.7847
$ws.Range("L138").Value = 14813.3541
$ws.Range("N138").Value = -25093.3541

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6334.0835
$ws.Range("I32").Value = 6334.0835
$ws.Range("K32").Value = 6334.0835
$ws.Range("M32").Value = -6047.0835

$ws.Range("H45").Value = 913.125
$ws.Range("I45").Value = 869
$ws.Range("J45").Value = 957.25
$ws.Range("K45").Value = 869
$ws.Range("L45").Value = 957.25
$ws.Range("M45").Value = -492
$ws.Range("N45").Value = -1711.25

$ws.Range("H61").Value = 4132
$ws.Range("I61").Value = 2545.2727
$ws.Range("K61").Value = 2545.2727
$ws.Range("M61").Value = -2333.2727

$ws.Range("H74").Value = 72023.11
$ws.Range("I74").Value = 74579.516
$ws.Range("K74").Value = 74579.516
$ws.Range("M74").Value = -73705.516

$ws.Range("H77").Value = 72023.11
$ws.Range("I77").Value = 74579.516
$ws.Range("K77").Value = 372897.58
$ws.Range("M77").Value = -368529.58

$ws.Range("H102").Value = 1898.9395
$ws.Range("I102").Value = 1523.44
$ws.Range("K102").Value = 1523.44
$ws.Range("M102").Value = 98.55999999999995

$ws.Range("H110").Value = 9552.541999999999
$ws.Range("I110").Value = 9995.933999999999
$ws.Range("J110").Value = 8813.556
$ws.Range("K110").Value = 9995.933999999999
$ws.Range("L110").Value = 8813.556
$ws.Range("M110").Value = -7950.933999999999
$ws.Range("N110").Value = -12903.556

$ws.Range("H136").Value = 4132
$ws.Range("I136").Value = 2545.2727
$ws.Range("K136").Value = 7635.8181
$ws.Range("M136").Value = -5085.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1873.4
$ws.Range("I64").Value = 586.5
$ws.Range("J64").Value = 2731.3333
$ws.Range("K64").Value = 586.5
$ws.Range("L64").Value = 2731.3333
$ws.Range("M64").Value = -361.5
$ws.Range("N64").Value = -3181.3333

$ws.Range("H67").Value = 1873.4
$ws.Range("I67").Value = 586.5
$ws.Range("J67").Value = 2731.3333
$ws.Range("K67").Value = 586.5
$ws.Range("L67").Value = 2731.3333
$ws.Range("M67").Value = 193.5
$ws.Range("N67").Value = -4291.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2223.6897
$ws.Range("I16").Value = 2223.6897
$ws.Range("K16").Value = 2223.6897
$ws.Range("M16").Value = -1936.6897

$ws.Range("H22").Value = 527.53845
$ws.Range("I22").Value = 388.25
$ws.Range("J22").Value = 2199
$ws.Range("K22").Value = 388.25
$ws.Range("L22").Value = 2199
$ws.Range("M22").Value = -38.25
$ws.Range("N22").Value = -2899

$ws.Range("H31").Value = 181515.97
$ws.Range("I31").Value = 224166.56
$ws.Range("J31").Value = 53564.2
$ws.Range("K31").Value = 224166.56
$ws.Range("L31").Value = 53564.2
$ws.Range("M31").Value = -223871.56
$ws.Range("N31").Value = -54154.2

$ws.Range("H34").Value = 181515.97
$ws.Range("I34").Value = 224166.56
$ws.Range("J34").Value = 53564.2
$ws.Range("K34").Value = 224166.56
$ws.Range("L34").Value = 53564.2
$ws.Range("M34").Value = -223964.56
$ws.Range("N34").Value = -53968.2

$ws.Range("H94").Value = 1611.8182
$ws.Range("I94").Value = 1404.25
$ws.Range("J94").Value = 1730.4286
$ws.Range("K94").Value = 1404.25
$ws.Range("L94").Value = 1730.4286
$ws.Range("M94").Value = -953.25
$ws.Range("N94").Value = -2632.4286

$ws.Range("H107").Value = 6315.84
$ws.Range("I107").Value = 1166.6666
$ws.Range("K107").Value = 1166.6666
$ws.Range("M107").Value = 753.3334

$ws.Range("H113").Value = 2223.6897
$ws.Range("I113").Value = 2223.6897
$ws.Range("K113").Value = 2223.6897
$ws.Range("M113").Value = -53.6896999999999

$ws.Range("H122").Value = 897
$ws.Range("I122").Value = 898.5
$ws.Range("J122").Value = 891
$ws.Range("K122").Value = 2695.5
$ws.Range("L122").Value = 2673
$ws.Range("M122").Value = -245.5
$ws.Range("N122").Value = -7573

$ws.Range("H132").Value = 2165.647
$ws.Range("I132").Value = 2196.9688
$ws.Range("J132").Value = 1664.5
$ws.Range("K132").Value = 6590.9064
$ws.Range("L132").Value = 4993.5
$ws.Range("M132").Value = -4060.9064
$ws.Range("N132").Value = -10053.5

$ws.Range("H141").Value = 112570.5
$ws.Range("J141").Value = 112570.5
$ws.Range("L141").Value = 112570.5
$ws.Range("N141").Value = -122930.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 750.375
$ws.Range("I5").Value = 750.4286
$ws.Range("K5").Value = 2251.2858
$ws.Range("M5").Value = -2139.2858

$ws.Range("H18").Value = 789.4167
$ws.Range("I18").Value = 586.44446
$ws.Range("J18").Value = 1398.3334
$ws.Range("K18").Value = 1759.33338
$ws.Range("L18").Value = 4195.0002
$ws.Range("M18").Value = -1590.33338
$ws.Range("N18").Value = -4533.0002

$ws.Range("H44").Value = 2707.9546
$ws.Range("I44").Value = 1525
$ws.Range("J44").Value = 2970.8333
$ws.Range("K44").Value = 4575
$ws.Range("L44").Value = 8912.499899999999
$ws.Range("M44").Value = -4177
$ws.Range("N44").Value = -9708.499899999999

$ws.Range("H135").Value = 750.375
$ws.Range("I135").Value = 750.4286
$ws.Range("K135").Value = 6753.8574
$ws.Range("M135").Value = -4218.8574

$ws.Range("H138").Value = 3665.875
$ws.Range("I138").Value = 4698.5
$ws.Range("K138").Value = 14095.5
$ws.Range("M138").Value = -8955.5

$ws.Range("H139").Value = 2328.5386
$ws.Range("I139").Value = 1738.8572
$ws.Range("J139").Value = 3016.5
$ws.Range("K139").Value = 5216.571599999999
$ws.Range("L139").Value = 9049.5
$ws.Range("M139").Value = -76.57159999999931
$ws.Range("N139").Value = -19329.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9141.556
$ws.Range("J70").Value = 7749
$ws.Range("L70").Value = 7749
$ws.Range("N70").Value = -8289

$ws.Range("H73").Value = 9141.556
$ws.Range("J73").Value = 7749
$ws.Range("L73").Value = 7749
$ws.Range("N73").Value = -9621

$ws.Range("H96").Value = 46608
$ws.Range("J96").Value = 46608
$ws.Range("L96").Value = 46608
$ws.Range("N96").Value = -52100

$ws.Range("H122").Value = 2099.4138
$ws.Range("I122").Value = 2128.3076
$ws.Range("J122").Value = 1849
$ws.Range("K122").Value = 6384.9228
$ws.Range("L122").Value = 5547
$ws.Range("M122").Value = -3934.9228
$ws.Range("N122").Value = -10447

$ws.Range("H132").Value = 6631.31
$ws.Range("I132").Value = 6076.975
$ws.Range("J132").Value = 8848.65
$ws.Range("K132").Value = 18230.925
$ws.Range("L132").Value = 26545.95
$ws.Range("M132").Value = -15700.925
$ws.Range("N132").Value = -31605.95

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 25001176
$ws.Range("I100").Value = 25001176
$ws.Range("K100").Value = 25001176
$ws.Range("M100").Value = -25000635

$ws.Range("H132").Value = 8098.6387
$ws.Range("I132").Value = 8292.458000000001
$ws.Range("J132").Value = 7711
$ws.Range("K132").Value = 24877.374
$ws.Range("L132").Value = 23133
$ws.Range("M132").Value = -22347.374
$ws.Range("N132").Value = -28193

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 39698.4
$ws.Range("I49").Value = 39698.4
$ws.Range("K49").Value = 39698.4
$ws.Range("M49").Value = -39468.4

$ws.Range("H100").Value = 1085.4667
$ws.Range("I100").Value = 1072.8889
$ws.Range("K100").Value = 2145.7778
$ws.Range("M100").Value = -1604.7778

$ws.Range("H136").Value = 166874.61
$ws.Range("J136").Value = 4895
$ws.Range("L136").Value = 14685
$ws.Range("N136").Value = -19785
